# Slide 19 reads "homework:" + line-break + "Caesar cypher".
# Slide 20 is a near-duplicate reading "homework:" + line-break + "Caesar cypher questions".
# The older/shorter duplicate (slide 19) is removed, leaving the fuller
# "Caesar cypher questions" slide as the deck's new final slide.

$p = $ppt.ActivePresentation
$p.Slides.Item(19).Delete()
